$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 28-30 (columns A..T). F/G/H use $null to mean "empty string cell"
# (handled specially below, since the source rows have <c t="inlineStr"/> - present but empty).
$rows = @(
    @{
        row = 28
        A = 26; B = 3; C = 0; D = 0.003; E = "Regular"
        F = $null; G = $null; H = $null
        I = "<function relu at 0x10f4d69d8>"
        J = 0.9516000151634216
        K = 0.04560000076889992
        L = 0.003700000001117587
        M = 0.1825118958950043
        N = 6.906796932220459
        O = 0.04560000076889992
        P = "logs/results_282.log"
        Q = "weights/model_282.ckpt"
        R = "tb/282"
        S = "(6.9546156, 7.3834124, 9.092276, 9.433221, 9.5498905, 11.14911, 9.904368)"
        T = "(139.42224, 9.018682, 9.271418, 8.975029, 7.9992733, 7.393931, 7.053868, 10.015819)"
    },
    @{
        row = 29
        A = 27; B = 3; C = 0; D = 0.003; E = "Regular"
        F = $null; G = $null; H = $null
        I = "<function relu at 0x1100289d8>"
        J = 0.9430999755859375
        K = 0.06279999762773514
        L = 0.006200000178068876
        M = 0.2192680686712265
        N = 7.151318073272705
        O = 0.06279999762773514
        P = "logs/results_285.log"
        Q = "weights/model_285.ckpt"
        R = "tb/285"
        S = "(7.0936155, 7.6124697, 8.267413, 8.219525, 11.057663, 10.19839, 9.6345)"
        T = "(141.34113, 8.731318, 9.105043, 8.344593, 9.322138, 7.830576, 7.469233, 9.507704)"
    },
    @{
        row = 30
        A = 28; B = 3; C = 0; D = 0.003; E = "Regular"
        F = $null; G = $null; H = $null
        I = "<function relu at 0x116aa89d8>"
        J = 0.9474999904632568
        K = 0.07349999994039536
        L = 0.03050000034272671
        M = 0.2007102072238922
        N = 6.25114631652832
        O = 0.07349999994039536
        P = "logs/results_305.log"
        Q = "weights/model_305.ckpt"
        R = "tb/305"
        S = "(6.9461164, 7.606389, 7.896417, 8.789286, 9.165759, 8.689637, 8.321884)"
        T = "(138.02159, 8.655811, 9.216804, 9.208384, 7.499287, 7.888038, 7.191258, 13.024145)"
    }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $addr = "$col$($r.row)"
        $val = $r[$col]
        if ($val -eq $null) {
            # Write an empty-string text cell (not a truly blank cell): enter a lone
            # apostrophe (Excel's text-prefix escape) which commits as "", then strip
            # the quote-prefix formatting it implies so the cell keeps the default style.
            $ws.Range($addr).Value = "'"
            $ws.Range($addr).Style = "Normal"
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}
